$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("35C_10.0MPa").Name = "35C_10MPa"
$wb.Worksheets.Item("35C_20.1MPa").Name = "35C_20MPa"
$wb.Worksheets.Item("50C_10.1MPa").Name = "50C_10MPa"
$wb.Worksheets.Item("50C_20.1MPa").Name = "50C_20MPa"
